# Fill in the missing box-score / forecast data for the Fri, Feb 2, 2024
# slate (Sheet1 rows 105-114): away/home points, overtime flag, win/loss
# team, the model's forecast pick, whether the forecast was correct, and
# the |Away-Home| point-differential formula that already runs as a
# shared formula down column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 104 already carries the correctly-styled "Yes" forecast (its K cell
# uses the bold/centered style) - use it as the template we copy onto any
# new K cell whose forecast turned out correct.
$styleDonor = $ws.Range("K104")

$rows = @(
    @{ Row = 105; D = 136; F = 125; I = "Los Angeles Clippers";   J = "Detroit Pistons";         K = "Detroit Pistons";         Correct = $false },
    @{ Row = 106; D = 110; F = 102; I = "Miami Heat";             J = "Washington Wizards";       K = "Washington Wizards";       Correct = $false },
    @{ Row = 107; D = 120; F = 129; I = "Atlanta Hawks";          J = "Phoenix Suns";             K = "Atlanta Hawks";            Correct = $true  },
    @{ Row = 108; D = 133; F = 122; I = "Sacramento Kings";       J = "Indiana Pacers";           K = "Sacramento Kings";         Correct = $true  },
    @{ Row = 109; D = 106; F = 135; I = "Houston Rockets";        J = "Toronto Raptors";          K = "Toronto Raptors";          Correct = $false },
    @{ Row = 110; D = 121; F = 101; I = "Golden State Warriors";  J = "Memphis Grizzlies";        K = "Memphis Grizzlies";        Correct = $false },
    @{ Row = 111; D = 108; F = 106; I = "Orlando Magic";          J = "Minnesota Timberwolves";   K = "Orlando Magic";            Correct = $true  },
    @{ Row = 112; D = 106; F = 126; I = "Oklahoma City Thunder";  J = "Charlotte Hornets";        K = "Charlotte Hornets";        Correct = $false },
    @{ Row = 113; D = 114; F = 113; I = "New Orleans Pelicans";   J = "San Antonio Spurs";        K = "San Antonio Spurs";        Correct = $false },
    @{ Row = 114; D = 108; F = 120; I = "Denver Nuggets";         J = "Portland Trail Blazers";   K = "Portland Trail Blazers";   Correct = $false }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 4).Value = $r.D          # D: Away Pts
    $ws.Cells.Item($row, 6).Value = $r.F          # F: Home Pts
    $ws.Cells.Item($row, 7).Value = "NA"          # G: Overtime
    $ws.Cells.Item($row, 9).Value = $r.I           # I: Win
    $ws.Cells.Item($row, 10).Value = $r.J          # J: Loss

    $kCell = $ws.Cells.Item($row, 11)
    if ($r.Correct) {
        # Match the workbook's existing convention: a correct forecast's
        # K cell is bold + centered, so stamp that formatting first.
        $styleDonor.Copy()
        $kCell.PasteSpecial(-4122)
        $kCell.Value = $r.K
        $ws.Cells.Item($row, 12).Value = "Yes"     # L: Correct
    } else {
        $kCell.Value = $r.K
        $ws.Cells.Item($row, 12).Value = "No"      # L: Correct
    }
}

$excel.CutCopyMode = $false

# Column M's "|Away-Home| diff" formula is a shared formula currently
# spanning M67:M104 - extend it down across the newly-populated rows.
$ws.Range("M105:M114").Formula = "=ABS(D105-F105)"
